# Generate Report for Handoff
# Updates the localization-status workbook with a new handoff cycle:
#  - new source file guid (c149f494-704f-41c4-ac16-2033550e7afd.md)
#  - new handoff xliff names / timestamps
#  - handback info reset (blank file / default date)
#  - "Has metadata" flipped to True

$wb = $excel.ActiveWorkbook

$oldGuid = "d6a5534b-124a-4535-8f3c-cfb62ed6460e"
$newGuid = "c149f494-704f-41c4-ac16-2033550e7afd"

$oldFileName = "$oldGuid.md"
$newFileName = "$newGuid.md"
$oldPathName = "e2e\$oldGuid.md"
$newPathName = "e2e\$newGuid.md"

$newHoDate          = "2017-01-03 05:27:29"
$newZhHandoffXlf    = "$newGuid.af809464802b244e49d80b692ec4c2c1e9ce0171.zh-cn.xlf"
$newZhHandoffDate   = "2017-01-03 05:27:18"
$newDeHandoffXlf    = "$newGuid.af809464802b244e49d80b692ec4c2c1e9ce0171.de-de.xlf"
$resetHandbackDate  = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value2 = $newFileName
$wsOverview.Range("B2").Value2 = $newPathName
$wsOverview.Range("G2").Value2 = $newHoDate

foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Address() -eq "`$B`$2") {
        $h.TextToDisplay = $newPathName
    }
}

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value2 = $newFileName
$wsZh.Range("G2").Value2 = $newZhHandoffXlf
$wsZh.Range("H2").Value2 = $newZhHandoffDate

# Latest Target File (J2) - clear value and hyperlink, restore normal style
$wsZh.Range("J2").Style = "Normal"
$wsZh.Range("J2").Value2 = ""

# Latest Handback File (K2) - clear
$wsZh.Range("K2").Value2 = ""

# Latest Handback DateTime (L2) - reset to default date
$wsZh.Range("L2").Value2 = $resetHandbackDate

# Has metadata (Q2) False -> True (leading apostrophe forces text, not boolean)
$wsZh.Range("Q2").Value2 = "'True"

# Remove the Latest Target File hyperlink (J2), keep/update the A2 one
$zhToDelete = @()
foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Address() -eq "`$J`$2") {
        $zhToDelete += $h
    }
}
foreach ($d in $zhToDelete) { $d.Delete() }

foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Address() -eq "`$A`$2") {
        $h.TextToDisplay = $newFileName
    }
}

# Column widths for J/K (10/11) shrink from 40 to roughly 18.65 / 21.71
$wsZh.Columns.Item(10).ColumnWidth = 17.833333333333332
$wsZh.Columns.Item(11).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value2 = $newFileName
$wsDe.Range("G2").Value2 = $newDeHandoffXlf
$wsDe.Range("H2").Value2 = $newHoDate

# Latest Target File (J2) - clear value and hyperlink, restore normal style
$wsDe.Range("J2").Style = "Normal"
$wsDe.Range("J2").Value2 = ""

# Latest Handback File (K2) - clear
$wsDe.Range("K2").Value2 = ""

# Latest Handback DateTime (L2) - reset to default date
$wsDe.Range("L2").Value2 = $resetHandbackDate

# Has metadata (Q2) False -> True (leading apostrophe forces text, not boolean)
$wsDe.Range("Q2").Value2 = "'True"

# Remove the Latest Target File hyperlink (J2), keep/update the A2 one
$deToDelete = @()
foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Address() -eq "`$J`$2") {
        $deToDelete += $h
    }
}
foreach ($d in $deToDelete) { $d.Delete() }

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Address() -eq "`$A`$2") {
        $h.TextToDisplay = $newFileName
    }
}

# Column widths for J/K (10/11) shrink from 40 to roughly 18.65 / 21.71
$wsDe.Columns.Item(10).ColumnWidth = 17.833333333333332
$wsDe.Columns.Item(11).ColumnWidth = 20.833333333333332

"Done applying handoff report changes."
